# Generate Report for Handoff
# Inserts a new file entry (b10081ed-378b-439d-87ec-f0df2a4cd958.md) with
# status "Ready for handoff" right before the existing
# c7e89640-77b8-4370-970f-bf0463e5644d.md entry (which is also "Ready for
# handoff"), on all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Cells.Item(1,1).Value = "File Name"
$ws1.Cells.Item(1,2).Value = "zh-cn"
$ws1.Cells.Item(1,3).Value = "de-de"

$ws1.Cells.Item(2,1).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md"
$ws1.Cells.Item(2,2).Value = "In Translation"
$ws1.Cells.Item(2,3).Value = "In Translation"

$ws1.Cells.Item(3,1).Value = "ca5f0528-f586-4664-8880-ed70a9764e1e.md"
$ws1.Cells.Item(3,2).Value = "In Translation"
$ws1.Cells.Item(3,3).Value = "In Translation"

$ws1.Cells.Item(4,1).Value = "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md"
$ws1.Cells.Item(4,2).Value = "In Translation"
$ws1.Cells.Item(4,3).Value = "In Translation"

$ws1.Cells.Item(5,1).Value = "fec7be35-ac65-4cdb-a386-f39ce0811702.md"
$ws1.Cells.Item(5,2).Value = "In Translation"
$ws1.Cells.Item(5,3).Value = "In Translation"

# New row: b10081ed-378b-439d-87ec-f0df2a4cd958.md, ready for handoff
$ws1.Cells.Item(6,1).Value = "b10081ed-378b-439d-87ec-f0df2a4cd958.md"
$ws1.Cells.Item(6,2).Value = "Ready for handoff"
$ws1.Cells.Item(6,3).Value = "Ready for handoff"

# Existing row shifted down from 6 -> 7
$ws1.Cells.Item(7,1).Value = "c7e89640-77b8-4370-970f-bf0463e5644d.md"
$ws1.Cells.Item(7,2).Value = "Ready for handoff"
$ws1.Cells.Item(7,3).Value = "Ready for handoff"

# Existing row shifted down from 7 -> 8
$ws1.Cells.Item(8,1).Value = ".localization-config"
$ws1.Cells.Item(8,2).Value = "Not to be localized"
$ws1.Cells.Item(8,3).Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/e2e/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cded2d20919d64cac2a57e1e89d4fafc4699f2c9/e2e/ca5f0528-f586-4664-8880-ed70a9764e1e.md", "", "", "ca5f0528-f586-4664-8880-ed70a9764e1e.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md", "", "", "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/fec7be35-ac65-4cdb-a386-f39ce0811702.md", "", "", "fec7be35-ac65-4cdb-a386-f39ce0811702.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/47017fa415b99faa36d73e4359406cae44492fa5/e2e/b10081ed-378b-439d-87ec-f0df2a4cd958.md", "", "", "b10081ed-378b-439d-87ec-f0df2a4cd958.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/59b2ad617493011d33cca8c86234818ac60dcfd3/e2e/c7e89640-77b8-4370-970f-bf0463e5644d.md", "", "", "c7e89640-77b8-4370-970f-bf0463e5644d.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/.localization-config", "", "", ".localization-config") | Out-Null

$ws1.Range("A2:A8").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Cells.Item(1,1).Value = "Source File Name"
$ws2.Cells.Item(1,2).Value = "Status"
$ws2.Cells.Item(1,3).Value = "Latest Handoff File"
$ws2.Cells.Item(1,4).Value = "Latest Handoff Datetime"
$ws2.Cells.Item(1,5).Value = "Latest Target File"
$ws2.Cells.Item(1,6).Value = "Latest Handback File"
$ws2.Cells.Item(1,7).Value = "Latest Handback DateTime"
$ws2.Cells.Item(1,8).Value = "Handoff Reason"
$ws2.Cells.Item(1,9).Value = "Dependency From"

$ws2.Cells.Item(2,1).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md"
$ws2.Cells.Item(2,2).Value = "In Translation"
$ws2.Cells.Item(2,3).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf"
$ws2.Cells.Item(2,4).Value = "2016-01-29 02:05:27"
$ws2.Cells.Item(2,5).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md"
$ws2.Cells.Item(2,6).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf"
$ws2.Cells.Item(2,7).Value = "2016-01-29 02:06:46"
$ws2.Cells.Item(2,8).Value = "Include"

$ws2.Cells.Item(3,1).Value = "ca5f0528-f586-4664-8880-ed70a9764e1e.md"
$ws2.Cells.Item(3,2).Value = "In Translation"
$ws2.Cells.Item(3,3).Value = "ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.zh-cn.xlf"
$ws2.Cells.Item(3,4).Value = "2016-01-29 02:02:30"
$ws2.Cells.Item(3,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(3,8).Value = "Include"

$ws2.Cells.Item(4,1).Value = "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md"
$ws2.Cells.Item(4,2).Value = "In Translation"
$ws2.Cells.Item(4,3).Value = "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.zh-cn.xlf"
$ws2.Cells.Item(4,4).Value = "2016-01-29 02:00:36"
$ws2.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,8).Value = "Include"

$ws2.Cells.Item(5,1).Value = "fec7be35-ac65-4cdb-a386-f39ce0811702.md"
$ws2.Cells.Item(5,2).Value = "In Translation"
$ws2.Cells.Item(5,3).Value = "fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.zh-cn.xlf"
$ws2.Cells.Item(5,4).Value = "2016-01-29 02:00:36"
$ws2.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5,8).Value = "Include"

# New row: b10081ed-378b-439d-87ec-f0df2a4cd958.md, ready for handoff
$ws2.Cells.Item(6,1).Value = "b10081ed-378b-439d-87ec-f0df2a4cd958.md"
$ws2.Cells.Item(6,2).Value = "Ready for handoff"
$ws2.Cells.Item(6,3).Value = "b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.zh-cn.xlf"
$ws2.Cells.Item(6,4).Value = "2016-01-29 02:09:06"
$ws2.Cells.Item(6,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(6,8).Value = "Include"

# Existing row shifted down from 6 -> 7
$ws2.Cells.Item(7,1).Value = "c7e89640-77b8-4370-970f-bf0463e5644d.md"
$ws2.Cells.Item(7,2).Value = "Ready for handoff"
$ws2.Cells.Item(7,3).Value = "c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.zh-cn.xlf"
$ws2.Cells.Item(7,4).Value = "2016-01-29 02:03:13"
$ws2.Cells.Item(7,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(7,8).Value = "Include"

# Existing row shifted down from 7 -> 8
$ws2.Cells.Item(8,1).Value = ".localization-config"
$ws2.Cells.Item(8,2).Value = "Not to be localized"
$ws2.Cells.Item(8,4).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(8,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(8,8).Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/e2e/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df3c125202494c35fbe803a32696c77678f5adfc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fef8c773fed7da709965609c3e88bd86065f7098/e2e/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b1fd30633923aac2eb47ea1ed46f122c12c2f5cb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cded2d20919d64cac2a57e1e89d4fafc4699f2c9/e2e/ca5f0528-f586-4664-8880-ed70a9764e1e.md", "", "", "ca5f0528-f586-4664-8880-ed70a9764e1e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fdd0e51962edae5e27337e8763f11e7decc2931e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.zh-cn.xlf", "", "", "ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md", "", "", "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7159f3394f567e65aab429479b012be5fc739d91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.zh-cn.xlf", "", "", "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/fec7be35-ac65-4cdb-a386-f39ce0811702.md", "", "", "fec7be35-ac65-4cdb-a386-f39ce0811702.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7159f3394f567e65aab429479b012be5fc739d91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.zh-cn.xlf", "", "", "fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/47017fa415b99faa36d73e4359406cae44492fa5/e2e/b10081ed-378b-439d-87ec-f0df2a4cd958.md", "", "", "b10081ed-378b-439d-87ec-f0df2a4cd958.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47017fa415b99faa36d73e4359406cae44492fa5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.zh-cn.xlf", "", "", "b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/59b2ad617493011d33cca8c86234818ac60dcfd3/e2e/c7e89640-77b8-4370-970f-bf0463e5644d.md", "", "", "c7e89640-77b8-4370-970f-bf0463e5644d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c10f7811f98f9bfadbf946233d16a27c39e2b461/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.zh-cn.xlf", "", "", "c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/.localization-config", "", "", ".localization-config") | Out-Null

$ws2.Range("A2:A8").Style = "HyperLink"
$ws2.Range("C2:C7").Style = "HyperLink"
$ws2.Range("E2:E2").Style = "HyperLink"
$ws2.Range("F2:F2").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Cells.Item(1,1).Value = "Source File Name"
$ws3.Cells.Item(1,2).Value = "Status"
$ws3.Cells.Item(1,3).Value = "Latest Handoff File"
$ws3.Cells.Item(1,4).Value = "Latest Handoff Datetime"
$ws3.Cells.Item(1,5).Value = "Latest Target File"
$ws3.Cells.Item(1,6).Value = "Latest Handback File"
$ws3.Cells.Item(1,7).Value = "Latest Handback DateTime"
$ws3.Cells.Item(1,8).Value = "Handoff Reason"
$ws3.Cells.Item(1,9).Value = "Dependency From"

$ws3.Cells.Item(2,1).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md"
$ws3.Cells.Item(2,2).Value = "In Translation"
$ws3.Cells.Item(2,3).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf"
$ws3.Cells.Item(2,4).Value = "2016-01-29 02:05:43"
$ws3.Cells.Item(2,5).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md"
$ws3.Cells.Item(2,6).Value = "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf"
$ws3.Cells.Item(2,7).Value = "2016-01-29 02:07:14"
$ws3.Cells.Item(2,8).Value = "Include"

$ws3.Cells.Item(3,1).Value = "ca5f0528-f586-4664-8880-ed70a9764e1e.md"
$ws3.Cells.Item(3,2).Value = "In Translation"
$ws3.Cells.Item(3,3).Value = "ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.de-de.xlf"
$ws3.Cells.Item(3,4).Value = "2016-01-29 02:02:43"
$ws3.Cells.Item(3,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(3,8).Value = "Include"

$ws3.Cells.Item(4,1).Value = "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md"
$ws3.Cells.Item(4,2).Value = "In Translation"
$ws3.Cells.Item(4,3).Value = "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.de-de.xlf"
$ws3.Cells.Item(4,4).Value = "2016-01-29 02:01:08"
$ws3.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,8).Value = "Include"

$ws3.Cells.Item(5,1).Value = "fec7be35-ac65-4cdb-a386-f39ce0811702.md"
$ws3.Cells.Item(5,2).Value = "In Translation"
$ws3.Cells.Item(5,3).Value = "fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.de-de.xlf"
$ws3.Cells.Item(5,4).Value = "2016-01-29 02:01:08"
$ws3.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5,8).Value = "Include"

# New row: b10081ed-378b-439d-87ec-f0df2a4cd958.md, ready for handoff
$ws3.Cells.Item(6,1).Value = "b10081ed-378b-439d-87ec-f0df2a4cd958.md"
$ws3.Cells.Item(6,2).Value = "Ready for handoff"
$ws3.Cells.Item(6,3).Value = "b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.de-de.xlf"
$ws3.Cells.Item(6,4).Value = "2016-01-29 02:09:23"
$ws3.Cells.Item(6,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(6,8).Value = "Include"

# Existing row shifted down from 6 -> 7
$ws3.Cells.Item(7,1).Value = "c7e89640-77b8-4370-970f-bf0463e5644d.md"
$ws3.Cells.Item(7,2).Value = "Ready for handoff"
$ws3.Cells.Item(7,3).Value = "c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.de-de.xlf"
$ws3.Cells.Item(7,4).Value = "2016-01-29 02:03:27"
$ws3.Cells.Item(7,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(7,8).Value = "Include"

# Existing row shifted down from 7 -> 8
$ws3.Cells.Item(8,1).Value = ".localization-config"
$ws3.Cells.Item(8,2).Value = "Not to be localized"
$ws3.Cells.Item(8,4).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(8,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(8,8).Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/e2e/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5cd73997851fdeb9f18411f0fccd64529e33016/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fc39f6a6204f4c6f98952d6f72b1c86294837483/e2e/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a606ff923e91dfd827bed91a5e4f53968c445e19/ol-handback/OpenLocalizationTestOrg/oltest.de-de/tianzh/1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf", "", "", "1c296d45-5d47-4ea2-833a-6ce140ff6c3c.62f8b0e370056b5c60b509cc8f232cf323acc349.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cded2d20919d64cac2a57e1e89d4fafc4699f2c9/e2e/ca5f0528-f586-4664-8880-ed70a9764e1e.md", "", "", "ca5f0528-f586-4664-8880-ed70a9764e1e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cabfe249a6b83df2fd44477e2d1e31f67b8e9e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.de-de.xlf", "", "", "ca5f0528-f586-4664-8880-ed70a9764e1e.cf77e21e6246bcfe3044ed698610592054749e98.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md", "", "", "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3dd3b5701ad93428a25423275d4d16d7a63f5c39/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.de-de.xlf", "", "", "e19f98b5-3b95-4a47-b2a6-bcb720d6a288.ae8ff1cee695a8d272df219c042107f50ffacfcc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5504914dc9fe3c35c4a75e10bb8deac9f4647edf/e2e/fec7be35-ac65-4cdb-a386-f39ce0811702.md", "", "", "fec7be35-ac65-4cdb-a386-f39ce0811702.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3dd3b5701ad93428a25423275d4d16d7a63f5c39/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.de-de.xlf", "", "", "fec7be35-ac65-4cdb-a386-f39ce0811702.ce66bbc31b15f0e92f3836d660bf71bfbd19b0ae.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/47017fa415b99faa36d73e4359406cae44492fa5/e2e/b10081ed-378b-439d-87ec-f0df2a4cd958.md", "", "", "b10081ed-378b-439d-87ec-f0df2a4cd958.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47017fa415b99faa36d73e4359406cae44492fa5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.de-de.xlf", "", "", "b10081ed-378b-439d-87ec-f0df2a4cd958.47017fa415b99faa36d73e4359406cae44492fa5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/59b2ad617493011d33cca8c86234818ac60dcfd3/e2e/c7e89640-77b8-4370-970f-bf0463e5644d.md", "", "", "c7e89640-77b8-4370-970f-bf0463e5644d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaba2da6726b882121878de6fefad457f38bc098/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.de-de.xlf", "", "", "c7e89640-77b8-4370-970f-bf0463e5644d.2d3410e2e990df487a9d6bca381a555c8aba43cd.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/fc745b428994e5cf9f5e9828f8a07652d0e7541c/.localization-config", "", "", ".localization-config") | Out-Null

$ws3.Range("A2:A8").Style = "HyperLink"
$ws3.Range("C2:C7").Style = "HyperLink"
$ws3.Range("E2:E2").Style = "HyperLink"
$ws3.Range("F2:F2").Style = "HyperLink"

$wb.Save()
